$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.154.55"
$ws.Range("E2").Value = "  +4.51%  "
$ws.Range("D3").Value = "1.909.10"
$ws.Range("E3").Value = "  +5.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.82"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5098"
$ws.Range("E7").Value = "  +3.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.91"
$ws.Range("E8").Value = "  +3.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2960"
$ws.Range("E9").Value = "  +6.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06789"
$ws.Range("E10").Value = "  +5.84%  "
$ws.Range("D11").Value = "1.908.98"
$ws.Range("E11").Value = "  +5.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.23"
$ws.Range("E12").Value = "  +2.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07350"
$ws.Range("E13").Value = "  +3.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6907"
$ws.Range("E14").Value = "  +6.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.78"
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.875"
$ws.Range("E16").Value = "  +3.64%  "
$ws.Range("D17").Value = "30.162.27"
$ws.Range("E17").Value = "  +4.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008269"
$ws.Range("E18").Value = "  +11.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.0000"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.96"
$ws.Range("E20").Value = "  +5.81%  "
$ws.Range("D21").Value = "2.156.26"
$ws.Range("E21").Value = "  +5.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9988"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.815"
$ws.Range("E23").Value = "  +4.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.707"
$ws.Range("E24").Value = "  +6.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.131"
$ws.Range("E25").Value = "  +2.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.72"
$ws.Range("E26").Value = "  +2.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "135.75"
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.06"
$ws.Range("E28").Value = "  +2.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.996"
$ws.Range("E29").Value = "  +5.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.394"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.234"
$ws.Range("E31").Value = "  +1.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08808"
$ws.Range("E32").Value = "  +5.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.009"
$ws.Range("E33").Value = "  +4.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05060"
$ws.Range("E34").Value = "  +2.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.147"
$ws.Range("E35").Value = "  +4.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7168"
$ws.Range("E36").Value = "  +5.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.688"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.817"
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.280"
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9693"
$ws.Range("E40").Value = "  +1.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01688"
$ws.Range("E41").Value = "  +6.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.164"
$ws.Range("E42").Value = "  +1.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4291"
$ws.Range("E43").Value = "  +4.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.93"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9990"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.614"
$ws.Range("E46").Value = "  +5.74%  "
$ws.Range("E47").Value = "  +4.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05733"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.09"
$ws.Range("E49").Value = "  +4.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.395"
$ws.Range("E50").Value = "  +3.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3797"
$ws.Range("E51").Value = "  +4.60%  "
